# Update column F (dSF) values for specific rows to reflect
# the repulled data / push-all-data / mean-calculation changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    5  = -1
    12 = -1
    15 = -1
    18 = -1
    21 = 0
    23 = 1
    35 = 3
    37 = 1
    39 = 2
    42 = 0
    50 = 3
    52 = 3
    56 = 3
    58 = -1
    59 = 2
    65 = 0
    69 = -2
    75 = -10
    78 = -1
    80 = -1
    81 = 7
    83 = -3
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
